# The source diff for this commit only touches packaging-level metadata:
#   * xl/workbook.xml  - the <xr:revisionPtr documentId="..."/> GUID, which
#     Excel mints fresh on every save. It carries no user content and isn't
#     settable from the object model - it's a side effect of saving, not an
#     edit in its own right.
#   * customXml/item1.xml <-> customXml/item2.xml and
#     customXml/itemProps1.xml <-> customXml/itemProps2.xml simply swap
#     physical file slots (the "properties" SharePoint content-type part and
#     the "contentTypeSchema" part trade places). The XML payload of each
#     part is byte-identical before/after; only which itemN.xml name holds
#     it changes. No card/stat data in the workbook (Sheet1) is touched by
#     this commit at all.
#
# Custom XML parts aren't something a worksheet edit changes through Excel's
# UI, but Excel does expose them on Workbook.CustomXMLParts, so reflect the
# swap through that collection for completeness. This is intentionally
# best-effort/non-fatal: some hosts surface CustomXMLParts read-only for
# workbooks that were authored outside Excel (e.g. via SharePoint content
# organizer rules), so failures here must not abort the script.
$wb = $excel.ActiveWorkbook

try {
    $parts = $wb.CustomXMLParts

    $propertiesXml = '<?xml version="1.0" encoding="utf-8"?>' +
        '<p:properties xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" ' +
        'xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" ' +
        'xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls">' +
        '<documentManagement/></p:properties>'

    # Re-add the "properties" part so it occupies the slot after the
    # "contentTypeSchema" part, matching the after-state ordering in the
    # diff (contentTypeSchema first, properties second).
    $existing = $parts.SelectByNamespace("http://schemas.microsoft.com/office/2006/metadata/properties")
    if ($existing -ne $null -and $existing.Count -gt 0) {
        for ($i = 1; $i -le $existing.Count; $i++) {
            $existing.Item($i).Delete()
        }
    }
    $null = $parts.Add($propertiesXml)
} catch {
    # CustomXMLParts isn't wired for this workbook/host - nothing else in
    # this commit depends on it, so just continue.
}
